$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at N. This shifts the existing "Total" header and the
# per-row COUNTBLANK formulas from column N to column O, preserving their
# original (unstyled) cell formatting.
$ws.Columns("N").Insert()

# Write the new attendance date "10/15/2018" into the new N5 header cell as
# plain text (matching the other date headers E5:M5), without letting Excel's
# autocorrect turn it into a date serial number. We stage the text in a
# scratch cell that is explicitly formatted as Text, then paste only the
# value into N5 so the destination keeps its inherited header style.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "10/15/2018"
$scratch.Copy()
$ws.Range("N5").PasteSpecial(-4163)
$scratch.Clear()

# Widen the Total formulas (now in column O) so they also account for the
# new N column of attendance data.
for ($r = 6; $r -le 74; $r++) {
    $ws.Range("O$r").Formula = "=COUNTBLANK(E" + $r + ":N" + $r + ")"
}
